$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @(
    'Sát Phá Tham hội chiếu tại Quan Lộc',
    'Tử Phủ Vũ Tướng hội chiếu tại Quan Lộc',
    'Cơ Nguyệt Đồng Lương hội chiếu tại Quan Lộc',
    'Cự Nhật hội chiếu tại Quan Lộc',
    'Kình Đà hội chiếu tại Quan Lộc',
    'Xương Khúc hội chiếu tại Quan Lộc',
    'Hoả Linh hội chiếu tại Quan Lộc',
    'Không Kiếp hội chiếu tại Quan Lộc',
    'Quang Quý hội chiếu tại Quan Lộc',
    'Tả Hữu hội chiếu tại Quan Lộc',
    'Song Hao hội chiếu tại Quan Lộc',
    'Tang Hổ hội chiếu tại Quan Lộc',
    'Khốc Hư hội chiếu tại Quan Lộc',
    'Hình Riêu hội chiếu tại Quan Lộc',
    'Thai Toạ hội chiếu tại Quan Lộc',
    'Đào Hồng hội chiếu tại Quan Lộc',
    'Ấn Phù hội chiếu tại Quan Lộc',
    'Song Hao Quyền Lộc Kiếp Hoả hội chiếu tại Quan Lộc',
    'Tử Phủ Vũ Tướng Xương Khúc Khôi Việt Tả Hữu Khoa Quyền Lộc Long hội chiếu tại Quan Lộc',
    'Tử Phủ Vũ Tướng Tả Hữu Khoa Quyền Lộc Long Phượng hội chiếu tại Quan Lộc',
    'Tử Khúc Phá Dương Đà hội chiếu tại Quan Lộc',
    'Cơ Nguyệt Đồng Lương gặp Xương Khúc Tả Hữu hội chiếu tại Quan Lộc',
    'Cơ Nguyệt Đồng Lương Khoa Tả Hữu Quang Quý Quan Phúc hội chiếu tại Quan Lộc',
    'Sát Quyền hội chiếu tại Quan Lộc',
    'Lộc Mã hội chiếu tại Quan Lộc',
    'Kiếp Hư Hao Quyền hội chiếu tại Quan Lộc',
    'Tuế Hổ Phù Xương Khúc hội chiếu tại Quan Lộc',
    'Xương Khúc Tấu Long Phượng hội chiếu tại Quan Lộc',
    'Đào Hồng Riêu Tấu Cơ Vũ hội chiếu tại Quan Lộc',
    'Binh Hình Tướng Ấn hội chiếu tại Quan Lộc',
    'Hổ Tấu hội chiếu tại Quan Lộc',
    'Hình Riêu Y hội chiếu tại Quan Lộc',
    'Mã Hỏa Linh hội chiếu tại Quan Lộc',
    'Thai Tọa Hồng Đào hội chiếu tại Quan Lộc',
    'Tả Hữu Không Kiếp hội chiếu tại Quan Lộc',
    'Tả Hữu Binh Tướng hội chiếu tại Quan Lộc',
    'Đào Quyền hội chiếu tại Quan Lộc',
    'Đào Hồng Tả Cái Hữu hội chiếu tại Quan Lộc',
    'Quan Phúc Quang Tấu hội chiếu tại Quan Lộc',
    'Đào Hồng Xương Khúc hội chiếu tại Quan Lộc',
    'Đào Tử Phủ hội chiếu tại Quan Lộc',
    'Hổ Kình Sát hội chiếu tại Quan Lộc',
    'Hổ Tang Kiếp hội chiếu tại Quan Lộc',
    'Hổ Tang Không Kiếp hội chiếu tại Quan Lộc',
    'Âm Dương Lương hội chiếu tại Quan Lộc',
    'Cơ Lương Gia Hội hội chiếu tại Quan Lộc',
    'Nhật Chiếu Lôi Môn hội chiếu tại Quan Lộc',
    'Tả Hữu Xương Khúc hội chiếu tại Quan Lộc',
    'Tham Linh Triều Viên hội chiếu tại Quan Lộc',
    'Cự Hỏa Linh hội chiếu tại Quan Lộc',
    'Hồng Đào Kỵ hội chiếu tại Quan Lộc',
    'Cự Đồng Hình hội chiếu tại Quan Lộc',
    'Lương Phá hội chiếu tại Quan Lộc',
    'Khôi Việt hội chiếu tại Quan Lộc',
    'Tham Vũ Hỏa hội chiếu tại Quan Lộc',
    'Cự Kỵ hội chiếu tại Quan Lộc',
    'Kình Đà Hỏa hội chiếu tại Quan Lộc',
    'Kình Đà Không Kiếp hội chiếu tại Quan Lộc',
    'Đào Hồng Hỷ hội chiếu tại Quan Lộc',
    'Đào Hồng Kiếp Sát hội chiếu tại Quan Lộc',
    'Đào Hồng Kỵ hội chiếu tại Quan Lộc',
    'Đào Hồng Tả Phù Hữu Bật hội chiếu tại Quan Lộc',
    'Đào Hồng Tả Phù Hữu Bật Khoa Quyền Lộc hội chiếu tại Quan Lộc',
    'Đào Hồng Tả Phù Hữu Bật Khoa Quyền Lộc Long Trì Phượng Các hội chiếu tại Quan Lộc',
    'Lưu Hà Kiếp Sát hội chiếu tại Quan Lộc',
    'Phục Không Kiếp hội chiếu tại Quan Lộc',
    'Xương Khúc Khôi Việt hội chiếu tại Quan Lộc',
    'Khoa Quyền hội chiếu tại Quan Lộc',
    'Lộc Quyền hội chiếu tại Quan Lộc',
    'Tướng Binh Đào Hồng hội chiếu tại Quan Lộc',
    'Lương Khốc Tuế hội chiếu tại Quan Lộc',
    'Tướng Binh hội chiếu tại Quan Lộc',
    'Cự Tang hội chiếu tại Quan Lộc',
    'Cự Hỏa hội chiếu tại Quan Lộc',
    'Phá Hình Kỵ hội chiếu tại Quan Lộc',
    'Tang Trực Tuế hội chiếu tại Quan Lộc',
    'Tuế Xương Khúc hội chiếu tại Quan Lộc',
    'Tử Vi Tả Hữu hội chiếu tại Quan Lộc',
    'Tử Vi Tang Tả Hữu hội chiếu tại Quan Lộc'
)

$startRow = 4209
for ($i = 0; $i -lt $newValues.Count; $i++) {
    $row = $startRow + $i
    $val = $newValues[$i]
    $ws.Cells.Item($row, 1).Value = $val
    $ws.Cells.Item($row, 2).Value = $val
}

$ws.Range("B4209:B4287").Select()
